$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4: change from text "123456789" to the numeric value 123456789
$ws.Range("B4").Value = 123456789

# Append new row 5
$ws.Range("A5").Value = "Venkateswarrao V"

# B5 must stay text "99999" (not auto-converted to a number)
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "99999"
$ws.Range("B5").Style = "Normal"

# C5 is numeric 14
$ws.Range("C5").Value = 14

# D5 must stay text "2025-08-16" (not auto-converted to a date serial)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2025-08-16"
$ws.Range("D5").Style = "Normal"
